# SensorsSchedulerOutLight.xlsx - "OutAssignment" sheet update
# (recent changes - prior to the presentation and BFS-like configuration generation)
#
# The sheet grows from 20 data rows (A3:E22) to 27 data rows (A3:E29):
# existing rows are renumbered/recomputed and 7 new rows are appended.
# Rewrite the whole A3:E29 block in one shot, then restore the saved
# selection (A3:E13) that the author left active on this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OutAssignment")

$rows = @(
    @(2,3,2,6,14),
    @(3,8,4,9.75,10.75),
    @(3,0,0,10.75,11.006944444444445),
    @(3,2,1,11.006944444444445,19.006944444444443),
    @(4,4,2,12,20),
    @(13,16,6,15,16),
    @(13,0,0,16,18.5),
    @(13,15,5,18.5,19.75),
    @(5,19,3,18.25,19.25),
    @(5,0,0,19.25,19.354166666666668),
    @(5,11,3,19.354166666666668,20.354166666666668),
    @(2,10,5,12.645833333333334,13.645833333333334),
    @(4,4,2,12,20),
    @(13,13,5,15,16.5),
    @(13,0,0,16.5,16.604166666666668),
    @(13,16,6,16.604166666666668,17.604166666666668),
    @(13,0,0,17.604166666666668,18),
    @(13,19,5,18,19),
    @(13,0,0,19,19.076388888888889),
    @(13,12,5,19.076388888888889,20.576388888888889),
    @(5,14,3,18.25,19.25),
    @(5,0,0,19.25,19.604166666666668),
    @(5,11,3,19.604166666666668,20.604166666666668),
    @(13,16,6,15,16),
    @(13,0,0,16,18),
    @(13,14,5,18,19),
    @(5,11,3,19,20)
)

$numRows = $rows.Count
$numCols = 5

$data = New-Object 'object[,]' $numRows, $numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $data[$i, $j] = $rows[$i][$j]
    }
}

$firstDataRow = 3
$lastDataRow = $firstDataRow + $numRows - 1
$ws.Range("A$firstDataRow`:E$lastDataRow").Value = $data

# Restore the sheet's saved selection.
$ws.Range("A3:E13").Select() | Out-Null
